$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the affected region first so we can rewrite values freely
$ws.Range("A94:E149").UnMerge() | Out-Null

# Write new cell values for rows 94-149 (block reorder: HINTERLANDS, DARK_AGES, GUILDS, ADVENTURES, CORNUCOPIA_GUILDS)
$ws.Range("A94").Value = "HINTERLANDS"
$ws.Range("B94").Value = "1ST"
$ws.Range("C94").Value = "TWO"
$ws.Range("D94").Value = "ACTION"
$ws.Range("E94").Value = 1
$ws.Range("A95").Value = ""
$ws.Range("B95").Value = ""
$ws.Range("C95").Value = ""
$ws.Range("D95").Value = "TREASURE"
$ws.Range("E95").Value = 1
$ws.Range("A96").Value = ""
$ws.Range("B96").Value = ""
$ws.Range("C96").Value = "THREE"
$ws.Range("D96").Value = ""
$ws.Range("E96").Value = 1
$ws.Range("A97").Value = ""
$ws.Range("B97").Value = ""
$ws.Range("C97").Value = ""
$ws.Range("D97").Value = "ACTION"
$ws.Range("E97").Value = 3
$ws.Range("A98").Value = ""
$ws.Range("B98").Value = ""
$ws.Range("C98").Value = "FOUR"
$ws.Range("D98").Value = "ACTION"
$ws.Range("E98").Value = 3
$ws.Range("A99").Value = ""
$ws.Range("B99").Value = ""
$ws.Range("C99").Value = "FIVE"
$ws.Range("D99").Value = "ACTION"
$ws.Range("E99").Value = 5
$ws.Range("A100").Value = ""
$ws.Range("B100").Value = ""
$ws.Range("C100").Value = ""
$ws.Range("D100").Value = "ACTION - ATTACK"
$ws.Range("E100").Value = 1
$ws.Range("A101").Value = ""
$ws.Range("B101").Value = ""
$ws.Range("C101").Value = "OTHER"
$ws.Range("D101").Value = ""
$ws.Range("E101").Value = 1
$ws.Range("A102").Value = ""
$ws.Range("B102").Value = ""
$ws.Range("C102").Value = ""
$ws.Range("D102").Value = "ACTION"
$ws.Range("E102").Value = 1
$ws.Range("A103").Value = ""
$ws.Range("B103").Value = "2ND"
$ws.Range("C103").Value = "THREE"
$ws.Range("D103").Value = "ACTION"
$ws.Range("E103").Value = 1
$ws.Range("A104").Value = ""
$ws.Range("B104").Value = ""
$ws.Range("C104").Value = "FOUR"
$ws.Range("D104").Value = "ACTION"
$ws.Range("E104").Value = 3
$ws.Range("A105").Value = ""
$ws.Range("B105").Value = ""
$ws.Range("C105").Value = "FIVE"
$ws.Range("D105").Value = "ACTION"
$ws.Range("E105").Value = 2
$ws.Range("A106").Value = ""
$ws.Range("B106").Value = ""
$ws.Range("C106").Value = ""
$ws.Range("D106").Value = "ACTION - ATTACK"
$ws.Range("E106").Value = 2
$ws.Range("A107").Value = ""
$ws.Range("B107").Value = ""
$ws.Range("C107").Value = ""
$ws.Range("D107").Value = "ATTACK - TREASURE"
$ws.Range("E107").Value = 1
$ws.Range("A108").Value = ""
$ws.Range("B108").Value = "1RC"
$ws.Range("C108").Value = "TWO"
$ws.Range("D108").Value = "ACTION"
$ws.Range("E108").Value = 1
$ws.Range("A109").Value = ""
$ws.Range("B109").Value = ""
$ws.Range("C109").Value = "THREE"
$ws.Range("D109").Value = "ACTION - ATTACK"
$ws.Range("E109").Value = 1
$ws.Range("A110").Value = ""
$ws.Range("B110").Value = ""
$ws.Range("C110").Value = "FOUR"
$ws.Range("D110").Value = ""
$ws.Range("E110").Value = 1
$ws.Range("A111").Value = ""
$ws.Range("B111").Value = ""
$ws.Range("C111").Value = ""
$ws.Range("D111").Value = "ACTION"
$ws.Range("E111").Value = 1
$ws.Range("A112").Value = ""
$ws.Range("B112").Value = ""
$ws.Range("C112").Value = ""
$ws.Range("D112").Value = "ACTION - ATTACK"
$ws.Range("E112").Value = 1
$ws.Range("A113").Value = ""
$ws.Range("B113").Value = ""
$ws.Range("C113").Value = "FIVE"
$ws.Range("D113").Value = "ACTION"
$ws.Range("E113").Value = 2
$ws.Range("A114").Value = ""
$ws.Range("B114").Value = ""
$ws.Range("C114").Value = ""
$ws.Range("D114").Value = "TREASURE"
$ws.Range("E114").Value = 2
$ws.Range("A115").Value = "DARK_AGES"
$ws.Range("B115").Value = "1ST"
$ws.Range("C115").Value = "TWO"
$ws.Range("D115").Value = "ACTION"
$ws.Range("E115").Value = 3
$ws.Range("A116").Value = ""
$ws.Range("B116").Value = ""
$ws.Range("C116").Value = "THREE"
$ws.Range("D116").Value = "ACTION"
$ws.Range("E116").Value = 5
$ws.Range("A117").Value = ""
$ws.Range("B117").Value = ""
$ws.Range("C117").Value = ""
$ws.Range("D117").Value = "ACTION - ATTACK"
$ws.Range("E117").Value = 1
$ws.Range("A118").Value = ""
$ws.Range("B118").Value = ""
$ws.Range("C118").Value = "FOUR"
$ws.Range("D118").Value = ""
$ws.Range("E118").Value = 1
$ws.Range("A119").Value = ""
$ws.Range("B119").Value = ""
$ws.Range("C119").Value = ""
$ws.Range("D119").Value = "ACTION"
$ws.Range("E119").Value = 7
$ws.Range("A120").Value = ""
$ws.Range("B120").Value = ""
$ws.Range("C120").Value = ""
$ws.Range("D120").Value = "ACTION - ATTACK - LOOTER"
$ws.Range("E120").Value = 1
$ws.Range("A121").Value = ""
$ws.Range("B121").Value = ""
$ws.Range("C121").Value = ""
$ws.Range("D121").Value = "ACTION - LOOTER"
$ws.Range("E121").Value = 1
$ws.Range("A122").Value = ""
$ws.Range("B122").Value = ""
$ws.Range("C122").Value = "FIVE"
$ws.Range("D122").Value = "ACTION"
$ws.Range("E122").Value = 8
$ws.Range("A123").Value = ""
$ws.Range("B123").Value = ""
$ws.Range("C123").Value = ""
$ws.Range("D123").Value = "ACTION - ATTACK"
$ws.Range("E123").Value = 2
$ws.Range("A124").Value = ""
$ws.Range("B124").Value = ""
$ws.Range("C124").Value = ""
$ws.Range("D124").Value = "ACTION - ATTACK - LOOTER"
$ws.Range("E124").Value = 1
$ws.Range("A125").Value = ""
$ws.Range("B125").Value = ""
$ws.Range("C125").Value = ""
$ws.Range("D125").Value = "TREASURE"
$ws.Range("E125").Value = 1
$ws.Range("A126").Value = ""
$ws.Range("B126").Value = ""
$ws.Range("C126").Value = "OTHER"
$ws.Range("D126").Value = "ACTION"
$ws.Range("E126").Value = 3
$ws.Range("A127").Value = "GUILDS"
$ws.Range("B127").Value = "1ST"
$ws.Range("C127").Value = "TWO"
$ws.Range("D127").Value = "ACTION"
$ws.Range("E127").Value = 2
$ws.Range("A128").Value = ""
$ws.Range("B128").Value = ""
$ws.Range("C128").Value = "FOUR"
$ws.Range("D128").Value = "ACTION"
$ws.Range("E128").Value = 3
$ws.Range("A129").Value = ""
$ws.Range("B129").Value = ""
$ws.Range("C129").Value = "FIVE"
$ws.Range("D129").Value = "ACTION"
$ws.Range("E129").Value = 4
$ws.Range("A130").Value = ""
$ws.Range("B130").Value = ""
$ws.Range("C130").Value = ""
$ws.Range("D130").Value = "ACTION - ATTACK"
$ws.Range("E130").Value = 1
$ws.Range("A131").Value = ""
$ws.Range("B131").Value = "1RC"
$ws.Range("C131").Value = "THREE"
$ws.Range("D131").Value = "ACTION"
$ws.Range("E131").Value = 1
$ws.Range("A132").Value = ""
$ws.Range("B132").Value = ""
$ws.Range("C132").Value = ""
$ws.Range("D132").Value = "TREASURE"
$ws.Range("E132").Value = 1
$ws.Range("A133").Value = ""
$ws.Range("B133").Value = ""
$ws.Range("C133").Value = "FOUR"
$ws.Range("D133").Value = "ACTION - ATTACK"
$ws.Range("E133").Value = 1
$ws.Range("A134").Value = "ADVENTURES"
$ws.Range("B134").Value = "1ST"
$ws.Range("C134").Value = "TWO"
$ws.Range("D134").Value = "ACTION"
$ws.Range("E134").Value = 4
$ws.Range("A135").Value = ""
$ws.Range("B135").Value = ""
$ws.Range("C135").Value = ""
$ws.Range("D135").Value = "TREASURE"
$ws.Range("E135").Value = 1
$ws.Range("A136").Value = ""
$ws.Range("B136").Value = ""
$ws.Range("C136").Value = "THREE"
$ws.Range("D136").Value = "ACTION"
$ws.Range("E136").Value = 1
$ws.Range("A137").Value = ""
$ws.Range("B137").Value = ""
$ws.Range("C137").Value = ""
$ws.Range("D137").Value = "ACTION - DURATION"
$ws.Range("E137").Value = 4
$ws.Range("A138").Value = ""
$ws.Range("B138").Value = ""
$ws.Range("C138").Value = "FOUR"
$ws.Range("D138").Value = "ACTION"
$ws.Range("E138").Value = 7
$ws.Range("A139").Value = ""
$ws.Range("B139").Value = ""
$ws.Range("C139").Value = "FIVE"
$ws.Range("D139").Value = "ACTION"
$ws.Range("E139").Value = 6
$ws.Range("A140").Value = ""
$ws.Range("B140").Value = ""
$ws.Range("C140").Value = ""
$ws.Range("D140").Value = "ACTION - ATTACK"
$ws.Range("E140").Value = 1
$ws.Range("A141").Value = ""
$ws.Range("B141").Value = ""
$ws.Range("C141").Value = ""
$ws.Range("D141").Value = "ACTION - ATTACK - DURATION"
$ws.Range("E141").Value = 3
$ws.Range("A142").Value = ""
$ws.Range("B142").Value = ""
$ws.Range("C142").Value = ""
$ws.Range("D142").Value = "ATTACK - TREASURE"
$ws.Range("E142").Value = 1
$ws.Range("A143").Value = ""
$ws.Range("B143").Value = ""
$ws.Range("C143").Value = ""
$ws.Range("D143").Value = "TREASURE"
$ws.Range("E143").Value = 1
$ws.Range("A144").Value = ""
$ws.Range("B144").Value = ""
$ws.Range("C144").Value = "OTHER"
$ws.Range("D144").Value = "ACTION - DURATION"
$ws.Range("E144").Value = 1
$ws.Range("A145").Value = "CORNUCOPIA_GUILDS"
$ws.Range("B145").Value = "2ND"
$ws.Range("C145").Value = "TWO"
$ws.Range("D145").Value = "ACTION"
$ws.Range("E145").Value = 1
$ws.Range("A146").Value = ""
$ws.Range("B146").Value = ""
$ws.Range("C146").Value = "THREE"
$ws.Range("D146").Value = "ACTION"
$ws.Range("E146").Value = 2
$ws.Range("A147").Value = ""
$ws.Range("B147").Value = ""
$ws.Range("C147").Value = "FOUR"
$ws.Range("D147").Value = "ACTION"
$ws.Range("E147").Value = 1
$ws.Range("A148").Value = ""
$ws.Range("B148").Value = ""
$ws.Range("C148").Value = "FIVE"
$ws.Range("D148").Value = "ACTION"
$ws.Range("E148").Value = 2
$ws.Range("A149").Value = ""
$ws.Range("B149").Value = ""
$ws.Range("C149").Value = ""
$ws.Range("D149").Value = "ACTION - ATTACK"
$ws.Range("E149").Value = 2

# Re-apply merged cells for this region
$ws.Range("A94:A114").Merge() | Out-Null
$ws.Range("A115:A126").Merge() | Out-Null
$ws.Range("A127:A133").Merge() | Out-Null
$ws.Range("A134:A144").Merge() | Out-Null
$ws.Range("A145:A149").Merge() | Out-Null
$ws.Range("B94:B102").Merge() | Out-Null
$ws.Range("B103:B107").Merge() | Out-Null
$ws.Range("B108:B114").Merge() | Out-Null
$ws.Range("B115:B126").Merge() | Out-Null
$ws.Range("B127:B130").Merge() | Out-Null
$ws.Range("B131:B133").Merge() | Out-Null
$ws.Range("B134:B144").Merge() | Out-Null
$ws.Range("B145:B149").Merge() | Out-Null
$ws.Range("C94:C95").Merge() | Out-Null
$ws.Range("C96:C97").Merge() | Out-Null
$ws.Range("C99:C100").Merge() | Out-Null
$ws.Range("C101:C102").Merge() | Out-Null
$ws.Range("C105:C107").Merge() | Out-Null
$ws.Range("C110:C112").Merge() | Out-Null
$ws.Range("C113:C114").Merge() | Out-Null
$ws.Range("C116:C117").Merge() | Out-Null
$ws.Range("C118:C121").Merge() | Out-Null
$ws.Range("C122:C125").Merge() | Out-Null
$ws.Range("C129:C130").Merge() | Out-Null
$ws.Range("C131:C132").Merge() | Out-Null
$ws.Range("C134:C135").Merge() | Out-Null
$ws.Range("C136:C137").Merge() | Out-Null
$ws.Range("C139:C143").Merge() | Out-Null
$ws.Range("C148:C149").Merge() | Out-Null

Write-Output "done"